# Applies the 2024-06-01 data update to violent-crime-full-year.xlsx
# Generated from the canonical OOXML diff; updates K-column (2024) totals
# and a couple of I-column (2022) corrections across the Citywide Totals,
# By Neighborhood, and per-neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 3111
$ws.Range('K3').Value = 3088
$ws.Range('I4').Value = 1796
$ws.Range('K4').Value = 635
$ws.Range('K5').Value = 206
$ws.Range('K6').Value = 3653
$ws.Range('I7').Value = 26250
$ws.Range('K7').Value = 10693
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 85
$ws.Range('K5').Value = 18
$ws.Range('K7').Value = 312
$ws.Range('I8').Value = 1542
$ws.Range('K8').Value = 709
$ws.Range('K11').Value = 227
$ws.Range('K15').Value = 110
$ws.Range('K19').Value = 324
$ws.Range('K20').Value = 247
$ws.Range('K23').Value = 98
$ws.Range('K24').Value = 36
$ws.Range('K25').Value = 45
$ws.Range('K29').Value = 560
$ws.Range('K33').Value = 423
$ws.Range('K35').Value = 16
$ws.Range('K36').Value = 126
$ws.Range('K37').Value = 364
$ws.Range('K42').Value = 370
$ws.Range('K44').Value = 101
$ws.Range('K48').Value = 131
$ws.Range('K50').Value = 62
$ws.Range('K51').Value = 117
$ws.Range('K52').Value = 294
$ws.Range('K53').Value = 146
$ws.Range('K54').Value = 208
$ws.Range('K55').Value = 110
$ws.Range('K57').Value = 34
$ws.Range('K60').Value = 65
$ws.Range('K63').Value = 38
$ws.Range('K65').Value = 251
$ws.Range('K66').Value = 39
$ws.Range('K67').Value = 422
$ws.Range('K73').Value = 97
$ws.Range('K76').Value = 163
$ws.Range('K77').Value = 76
$ws.Range('K82').Value = 13
$ws.Range('K83').Value = 230
$ws.Range('K85').Value = 506
$ws.Range('K88').Value = 120
$ws.Range('K89').Value = 141
$ws.Range('K91').Value = 111
$ws.Range('K94').Value = 134
$ws.Range('K95').Value = 176
$ws.Range('K97').Value = 94
$ws.Range('K99').Value = 191
$ws.Range('I101').Value = 26250
$ws.Range('K101').Value = 10693
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K5').Value = 14
$ws.Range('K7').Value = 312
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K3').Value = 62
$ws.Range('K7').Value = 227
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K6').Value = 41
$ws.Range('K7').Value = 141
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 183
$ws.Range('K7').Value = 506
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K2').Value = 78
$ws.Range('K6').Value = 118
$ws.Range('K7').Value = 294
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K2').Value = 33
$ws.Range('K7').Value = 146
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 206
$ws.Range('K3').Value = 216
$ws.Range('I4').Value = 94
$ws.Range('I7').Value = 1542
$ws.Range('K7').Value = 709
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K3').Value = 72
$ws.Range('K7').Value = 230
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 118
$ws.Range('K3').Value = 153
$ws.Range('K6').Value = 121
$ws.Range('K7').Value = 423
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K6').Value = 45
$ws.Range('K7').Value = 176
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K3').Value = 125
$ws.Range('K7').Value = 364
$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 78
$ws.Range('K3').Value = 62
$ws.Range('K7').Value = 251
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K3').Value = 69
$ws.Range('K7').Value = 191
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 131
$ws.Range('K3').Value = 138
$ws.Range('K6').Value = 120
$ws.Range('K7').Value = 422
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K6').Value = 95
$ws.Range('K7').Value = 208
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 151
$ws.Range('K3').Value = 196
$ws.Range('K6').Value = 171
$ws.Range('K7').Value = 560
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K4').Value = 18
$ws.Range('K6').Value = 67
$ws.Range('K7').Value = 131
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K3').Value = 84
$ws.Range('K6').Value = 105
$ws.Range('K7').Value = 324
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K6').Value = 47
$ws.Range('K7').Value = 101
$ws = $wb.Worksheets.Item('River North')
$ws.Range('K6').Value = 95
$ws.Range('K7').Value = 163
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 94
$ws.Range('K6').Value = 141
$ws.Range('K7').Value = 370
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K4').Value = 5
$ws.Range('K7').Value = 110
$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('K4').Value = 3
$ws.Range('K7').Value = 36
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K2').Value = 35
$ws.Range('K3').Value = 30
$ws.Range('K7').Value = 98
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K3').Value = 52
$ws.Range('K6').Value = 25
$ws.Range('K7').Value = 111
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K3').Value = 69
$ws.Range('K7').Value = 247
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K2').Value = 50
$ws.Range('K7').Value = 126
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K4').Value = 13
$ws.Range('K7').Value = 134
$ws = $wb.Worksheets.Item('East Side')
$ws.Range('K6').Value = 7
$ws.Range('K7').Value = 45
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K2').Value = 38
$ws.Range('K7').Value = 110
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('K2').Value = 12
$ws.Range('K7').Value = 62
$ws = $wb.Worksheets.Item('North Center')
$ws.Range('K6').Value = 22
$ws.Range('K7').Value = 39
$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range('K6').Value = 12
$ws.Range('K7').Value = 16
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K3').Value = 24
$ws.Range('K7').Value = 97
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K3').Value = 24
$ws.Range('K4').Value = 6
$ws.Range('K7').Value = 85
$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K3').Value = 13
$ws.Range('K6').Value = 57
$ws.Range('K7').Value = 94
$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K6').Value = 60
$ws.Range('K7').Value = 120
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('K3').Value = 6
$ws.Range('K7').Value = 18
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K3').Value = 30
$ws.Range('K7').Value = 117
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('K6').Value = 21
$ws.Range('K7').Value = 34
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('K3').Value = 23
$ws.Range('K6').Value = 19
$ws.Range('K7').Value = 65
$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range('K3').Value = 3
$ws.Range('K6').Value = 13
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K2').Value = 38
$ws.Range('K7').Value = 76
